$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Table1")

# --- Row 11 ---
$ws.Range("A10:D10").Copy()
$row11 = $table.ListRows.Add()
$ws.Range("A11:D11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(11).RowHeight = 42.75
$row11.Range.Cells.Item(1, 1).Value = "The system will retrieve images of its workspace"
$row11.Range.Cells.Item(1, 2).Value = "The system will produce photographs within a specified work area"

# --- Row 12 ---
$ws.Range("A10:D10").Copy()
$row12 = $table.ListRows.Add()
$ws.Range("A12:D12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(12).RowHeight = 42.75
$row12.Range.Cells.Item(1, 1).Value = "The system will have a user interface"
$row12.Range.Cells.Item(1, 2).Value = "The system will accept user commands through a recognizable interface system"

$excel.CutCopyMode = $false

$ws.Range("F3").Select()
